# Weekly update: insert the newest price record at the top (row 2),
# pushing every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 2; existing rows 2..61 shift down to 3..62.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Match the date-cell number format used by the rest of column D.
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Populate the new record.
$ws.Cells.Item(2, 1).Value  = 1
$ws.Cells.Item(2, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value  = 44631
$ws.Cells.Item(2, 5).Value  = 15
$ws.Cells.Item(2, 6).Value  = 100112021
$ws.Cells.Item(2, 7).Value  = "Ají"
$ws.Cells.Item(2, 8).Value  = "Inferno"
$ws.Cells.Item(2, 9).Value  = "Primera"
$ws.Cells.Item(2, 10).Value = 140
$ws.Cells.Item(2, 11).Value = 19000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 19500
$ws.Cells.Item(2, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 1300
$ws.Cells.Item(2, 17).Value = 15
$ws.Cells.Item(2, 18).Value = "Hortaliza"
